$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: R007 -> R001 (Wohn Jick -> Jane Doe)
$ws.Range("A2").Value = "R001"
$ws.Range("B2").Value = "Jane Doe"
$ws.Range("C2").Value = "Graduate School"
$ws.Range("D2").Value = "Master"
$ws.Range("E2").Value = "Complete"
$ws.Range("F2").Value = "2023-05-04 23:21:13"
$ws.Range("G2").Value = 1

# Row 3: R008 -> R002 (Ran Dee -> John Doe)
$ws.Range("A3").Value = "R002"
$ws.Range("B3").Value = "John Doe"
$ws.Range("C3").Value = "Senior High School"
$ws.Range("D3").Value = "Accountancy, Business, and Management Strand"
$ws.Range("E3").Value = "Complete"
$ws.Range("F3").Value = "2023-05-04 23:21:34"
$ws.Range("G3").Value = 1

# Row 4: new row - R003 (Mark Doe)
$ws.Range("A4").Value = "R003"
$ws.Range("B4").Value = "Mark Doe"
$ws.Range("C4").Value = "College"
$ws.Range("D4").Value = "Bachelor of Science in Information Technology"
$ws.Range("E4").Value = "Complete"
$ws.Range("F4").Value = "2023-05-04 23:21:53"
$ws.Range("G4").Value = 1

# Row 5: new row - R004 (Son Doe)
$ws.Range("A5").Value = "R004"
$ws.Range("B5").Value = "Son Doe"
$ws.Range("C5").Value = "Junior High School"
$ws.Range("D5").Value = "Junior High School"
$ws.Range("E5").Value = "Complete"
$ws.Range("F5").Value = "2023-05-04 23:22:11"
$ws.Range("G5").Value = 1
